$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
# Row 8
$ws.Range("H8").Value = 42.4
$ws.Range("I8").Value = 42.153847
$ws.Range("J8").Value = 44
$ws.Range("K8").Value = 126.461541
$ws.Range("L8").Value = 132
$ws.Range("M8").Value = 12.538459
$ws.Range("N8").Value = -410
# Row 40
$ws.Range("H40").Value = 3314.2144
$ws.Range("I40").Value = 3278.8
$ws.Range("J40").Value = 3333.889
$ws.Range("K40").Value = 3278.8
$ws.Range("L40").Value = 3333.889
$ws.Range("M40").Value = -3103.8
$ws.Range("N40").Value = -3683.889
# Row 42
$ws.Range("H42").Value = 185.33333
$ws.Range("I42").Value = 185.33333
$ws.Range("K42").Value = 555.99999
$ws.Range("M42").Value = -325.99999
# Row 58
$ws.Range("H58").Value = 227.41667
$ws.Range("J58").Value = 69.8
$ws.Range("L58").Value = 209.4
$ws.Range("N58").Value = -509.4
# Row 98
$ws.Range("H98").Value = 2582.2307
$ws.Range("I98").Value = 2223.5
$ws.Range("K98").Value = 2223.5
$ws.Range("M98").Value = -725.5
# Row 107
$ws.Range("H107").Value = 696.1
$ws.Range("I107").Value = 711.8333
$ws.Range("J107").Value = 554.5
$ws.Range("K107").Value = 711.8333
$ws.Range("L107").Value = 554.5
$ws.Range("M107").Value = 1208.1667
$ws.Range("N107").Value = -4394.5
# Row 112
$ws.Range("H112").Value = 1311.5349
$ws.Range("J112").Value = 1274.0333
$ws.Range("L112").Value = 3822.0999
$ws.Range("N112").Value = -6038.0999
# Row 113
$ws.Range("H113").Value = 4174.778
$ws.Range("I113").Value = 3514.6
$ws.Range("K113").Value = 3514.6
$ws.Range("M113").Value = -260.5999999999999
# Row 122
$ws.Range("H122").Value = 2582.2307
$ws.Range("I122").Value = 2223.5
$ws.Range("K122").Value = 6670.5
$ws.Range("M122").Value = -4220.5
# Row 132
$ws.Range("H132").Value = 58713.453
$ws.Range("I132").Value = 64972.215
$ws.Range("K132").Value = 194916.645
$ws.Range("M132").Value = -192386.645
# Row 137
$ws.Range("H137").Value = 52633640
$ws.Range("I137").Value = 76925090
$ws.Range("K137").Value = 230775270
$ws.Range("M137").Value = -230772720

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1676
$ws.Range("I2").Value = 1732.5625
$ws.Range("J2").Value = 1546.7142
$ws.Range("K2").Value = 1732.5625
$ws.Range("L2").Value = 1546.7142
$ws.Range("M2").Value = -1619.5625
$ws.Range("N2").Value = -1772.7142
# Row 32
$ws.Range("H32").Value = 6099533.5
$ws.Range("I32").Value = 7693259.5
$ws.Range("K32").Value = 7693259.5
$ws.Range("M32").Value = -7692972.5
# Row 74
$ws.Range("H74").Value = 3790886.5
$ws.Range("I74").Value = 6946091
$ws.Range("K74").Value = 6946091
$ws.Range("M74").Value = -6945217
# Row 77
$ws.Range("H77").Value = 3790886.5
$ws.Range("I77").Value = 6946091
$ws.Range("K77").Value = 34730455
$ws.Range("M77").Value = -34726087
# Row 116
$ws.Range("H116").Value = 1676
$ws.Range("I116").Value = 1732.5625
$ws.Range("J116").Value = 1546.7142
$ws.Range("K116").Value = 1732.5625
$ws.Range("L116").Value = 1546.7142
$ws.Range("M116").Value = 561.4375
$ws.Range("N116").Value = -6134.7142
# Row 132
$ws.Range("H132").Value = 371217.9
$ws.Range("I132").Value = 432727.78
$ws.Range("J132").Value = 2158.6667
$ws.Range("K132").Value = 1298183.34
$ws.Range("L132").Value = 6476.000100000001
$ws.Range("M132").Value = -1295653.34
$ws.Range("N132").Value = -11536.0001

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1676
$ws.Range("I3").Value = 1732.5625
$ws.Range("J3").Value = 1546.7142
$ws.Range("K3").Value = 1732.5625
$ws.Range("L3").Value = 1546.7142
$ws.Range("M3").Value = -1618.5625
$ws.Range("N3").Value = -1774.7142
# Row 94
$ws.Range("H94").Value = 481
$ws.Range("I94").Value = 481
$ws.Range("K94").Value = 481
$ws.Range("M94").Value = -30
# Row 105
$ws.Range("H105").Value = 2000
$ws.Range("I105").Value = 2000
$ws.Range("K105").Value = 2000
$ws.Range("M105").Value = -253

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
# Row 4
$ws.Range("H4").Value = 19166.666
$ws.Range("J4").Value = 19166.666
$ws.Range("L4").Value = 19166.666
$ws.Range("N4").Value = -19390.666
# Row 58
$ws.Range("H58").Value = 9345092
$ws.Range("I58").Value = 15030998
$ws.Range("J58").Value = 3960.3572
$ws.Range("K58").Value = 15030998
$ws.Range("L58").Value = 3960.3572
$ws.Range("M58").Value = -15030795
$ws.Range("N58").Value = -4366.3572
# Row 62
$ws.Range("H62").Value = 4433
$ws.Range("I62").Value = 4649.5
$ws.Range("J62").Value = 4000
$ws.Range("K62").Value = 4649.5
$ws.Range("L62").Value = 4000
$ws.Range("M62").Value = -4025.5
$ws.Range("N62").Value = -5248
# Row 65
$ws.Range("H65").Value = 4433
$ws.Range("I65").Value = 4649.5
$ws.Range("J65").Value = 4000
$ws.Range("K65").Value = 23247.5
$ws.Range("L65").Value = 20000
$ws.Range("M65").Value = -20127.5
$ws.Range("N65").Value = -26240
# Row 107
$ws.Range("H107").Value = 939.125
$ws.Range("I107").Value = 804.1
$ws.Range("J107").Value = 1164.1666
$ws.Range("K107").Value = 804.1
$ws.Range("L107").Value = 1164.1666
$ws.Range("M107").Value = 1115.9
$ws.Range("N107").Value = -5004.1666
# Row 136
$ws.Range("H136").Value = 9345092
$ws.Range("I136").Value = 15030998
$ws.Range("J136").Value = 3960.3572
$ws.Range("K136").Value = 45092994
$ws.Range("L136").Value = 11881.0716
$ws.Range("M136").Value = -45090444
$ws.Range("N136").Value = -16981.0716

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("I4").Value = 147155550
$ws.Range("J4").Value = 18000056
$ws.Range("K4").Value = 441466650
$ws.Range("L4").Value = 54000168
$ws.Range("M4").Value = -441466538
$ws.Range("N4").Value = -54000392
# Row 55
$ws.Range("H55").Value = 1298.25
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()
# Row 57
$ws.Range("H57").Value = 2542.3333
$ws.Range("I57").Value = 1293
$ws.Range("K57").Value = 3879
$ws.Range("M57").Value = -3320
# Row 136
$ws.Range("H136").Value = 1399.4286
$ws.Range("I136").Value = 737.9231
$ws.Range("K136").Value = 2213.7693
$ws.Range("M136").Value = 2886.2307

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("H97").Value = 2242.5186
$ws.Range("I97").Value = 1259.4783
$ws.Range("K97").Value = 1259.4783
$ws.Range("M97").Value = -763.4783
# Row 99
$ws.Range("H99").Value = 19500
$ws.Range("I99").Value = 15000
$ws.Range("K99").Value = 15000
$ws.Range("M99").Value = -12754
# Row 113
$ws.Range("H113").Value = 2582.5715
$ws.Range("J113").Value = 3050.5
$ws.Range("L113").Value = 3050.5
$ws.Range("N113").Value = -7390.5
# Row 126
$ws.Range("H126").Value = 837568.1
$ws.Range("I126").Value = 1670686.2
$ws.Range("J126").Value = 4449.9
$ws.Range("K126").Value = 5012058.6
$ws.Range("L126").Value = 13349.7
$ws.Range("M126").Value = -5009588.6
$ws.Range("N126").Value = -18289.7

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 3834.182
$ws.Range("I7").Value = 3019.611
$ws.Range("J7").Value = 7499.75
$ws.Range("K7").Value = 3019.611
$ws.Range("L7").Value = 7499.75
$ws.Range("M7").Value = -2907.611
$ws.Range("N7").Value = -7723.75
# Row 24
$ws.Range("H24").Value = 1006
$ws.Range("I24").Value = 1006
$ws.Range("K24").Value = 1006
$ws.Range("M24").Value = -663
# Row 40
$ws.Range("H40").Value = 4125.2
$ws.Range("I40").Value = 3759.923
$ws.Range("K40").Value = 3759.923
$ws.Range("M40").Value = -3623.923
# Row 55
$ws.Range("H55").Value = 862.35297
$ws.Range("I55").Value = 287
$ws.Range("J55").Value = 1265.1
$ws.Range("K55").Value = 287
$ws.Range("L55").Value = 1265.1
$ws.Range("M55").Value = -114
$ws.Range("N55").Value = -1611.1
# Row 82
$ws.Range("H82").Value = 840
$ws.Range("I82").Value = 714.5714
$ws.Range("K82").Value = 714.5714
$ws.Range("M82").Value = -353.5714
# Row 85
$ws.Range("H85").Value = 840
$ws.Range("I85").Value = 714.5714
$ws.Range("K85").Value = 714.5714
$ws.Range("M85").Value = 533.4286
# Row 107
$ws.Range("H107").Value = 9494.333000000001
$ws.Range("I107").Value = 9494.333000000001
$ws.Range("K107").Value = 9494.333000000001
$ws.Range("M107").Value = -7574.333000000001
# Row 126
$ws.Range("H126").Value = 3834.182
$ws.Range("I126").Value = 3019.611
$ws.Range("J126").Value = 7499.75
$ws.Range("K126").Value = 9058.832999999999
$ws.Range("L126").Value = 22499.25
$ws.Range("M126").Value = -6588.832999999999
$ws.Range("N126").Value = -27439.25
# Row 132
$ws.Range("H132").Value = 20454862
$ws.Range("I132").Value = 21657854
$ws.Range("K132").Value = 64973562
$ws.Range("M132").Value = -64971032

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
# Row 126
$ws.Range("H126").Value = 2821.9092
$ws.Range("I126").Value = 2891.2104
$ws.Range("K126").Value = 8673.6312
$ws.Range("M126").Value = -6203.6312
# Row 136
$ws.Range("H136").Value = 10196463
$ws.Range("I136").Value = 11146650
$ws.Range("K136").Value = 33439950
$ws.Range("M136").Value = -33437400
